$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ps = $ws.PageSetup
Get-Member -InputObject $ps
